# "actualizacion de febrero hay un archivo mal"
# Update the SIPOT 3rd-trimester reporting sheet to the 4th-trimester period:
# move B8/C8 (period start/end) and J8/K8 (publish/update dates) forward by
# one quarter, and refresh the sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the sheet and scroll the window back towards the left/top (was
# parked at I2/I11 from the previous trimester's edit session).
$ws.Activate()
$win = $excel.Application.ActiveWindow
$win.ScrollRow    = 2
$win.ScrollColumn = 1

# New working cell for this edit session.
$ws.Range("C12").Select()

# Periodo que se informa: 07/01/2021-09/30/2021  ->  10/01/2021-12/31/2021
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561

# Fecha de validación / Fecha de actualización: 10/11/2021 -> 01/10/2022
$ws.Range("J8").Value = 44571
$ws.Range("K8").Value = 44571
